# The model-browser run populated the simulation table (C2:J35) with a
# uniform "representative" set of coefficients instead of the placeholder
# 1's / partially-tuned rows that were there before, then left the J column
# (dp_dwrbr_1) selected for review.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New uniform coefficient row used for every data row (rows 2-35),
# columns C..J (dp_dbnd_0, dp_dbnd_1, dp_dfrcw_0, dp_dfrcw_1,
# dp_dpcwr_0, dp_dpcwr_1, dp_dwrbr_0, dp_dwrbr_1).
$newValues = @(0.75, 0.75, 1.2, 1, 1.6, 1.1499999999999999, 1, 1)

$firstRow = 2
$lastRow = 35
$firstCol = 3   # column C

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($i = 0; $i -lt $newValues.Length; $i++) {
        $ws.Cells.Item($row, $firstCol + $i).Value = $newValues[$i]
    }
}

# Leave the dp_dwrbr_1 column (J2:J35) selected, as in the saved file.
$ws.Range("J2:J35").Select()
